$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 1.550905333333333
$ws.Range("N2").Value = 4.652716
$ws.Range("O2").Value = 0.3047927374587132
$ws.Range("P2").Value = 0.3047927374587132
$ws.Range("Q2").Value = 306.3710666146084
$ws.Range("R2").Value = 2757.339599531476
$ws.Range("S2").Value = 0.103285229593043
$ws.Range("T2").Value = 0.103285229593043

# Row 3
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("O3").Value = 0.3698654061072322
$ws.Range("P3").Value = 0.3698654061072323
$ws.Range("Q3").Value = 371.7807055303202
$ws.Range("R3").Value = 3346.026349772882
$ws.Range("S3").Value = 0.1253364292956105
$ws.Range("T3").Value = 0.1253364292956105

# Row 4
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("M4").Value = 0.260824
$ws.Range("N4").Value = 0.7824719999999999
$ws.Range("O4").Value = 0.05125861601369915
$ws.Range("P4").Value = 0.05125861601369915
$ws.Range("Q4").Value = 51.52405202382133
$ws.Range("R4").Value = 463.7164682143919
$ws.Range("S4").Value = 0.01737002648993138
$ws.Range("T4").Value = 0.01737002648993138

# Row 5
$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 1.162287333333333
$ws.Range("N5").Value = 3.486862
$ws.Range("O5").Value = 0.2284193176890152
$ws.Range("P5").Value = 0.2284193176890152
$ws.Range("Q5").Value = 229.6021571224091
$ws.Range("R5").Value = 2066.419414101682
$ws.Range("S5").Value = 0.07740454010716692
$ws.Range("T5").Value = 0.07740454010716692

# Row 6
$ws.Range("G6").Value = 197.5433703333333
$ws.Range("H6").Value = 592.6301109999999
$ws.Range("I6").Value = 0.3388703761585983
$ws.Range("J6").Value = 0.3388703761585982
$ws.Range("M6").Value = 0.232356
$ws.Range("N6").Value = 0.697068
$ws.Range("O6").Value = 0.04566392273134021
$ws.Range("P6").Value = 0.04566392273134022
$ws.Range("Q6").Value = 45.90038735717199
$ws.Range("R6").Value = 413.103486214548
$ws.Range("S6").Value = 0.01547415067284643
$ws.Range("T6").Value = 0.01547415067284642

# Row 7
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("M7").Value = 1.550905333333333
$ws.Range("N7").Value = 4.652716
$ws.Range("O7").Value = 0.3047927374587132
$ws.Range("P7").Value = 0.3047927374587132
$ws.Range("Q7").Value = 123.8045119060809
$ws.Range("R7").Value = 1114.240607154728
$ws.Range("S7").Value = 0.04173754910400692
$ws.Range("T7").Value = 0.04173754910400691

# Row 8
$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("O8").Value = 0.3698654061072322
$ws.Range("P8").Value = 0.3698654061072323
$ws.Range("S8").Value = 0.05064843630457296
$ws.Range("T8").Value = 0.05064843630457296

# Row 9
$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("M9").Value = 0.260824
$ws.Range("N9").Value = 0.7824719999999999
$ws.Range("O9").Value = 0.05125861601369915
$ws.Range("P9").Value = 0.05125861601369915
$ws.Range("Q9").Value = 20.82086334953066
$ws.Range("R9").Value = 187.387770145776
$ws.Range("S9").Value = 0.007019225657123819
$ws.Range("T9").Value = 0.007019225657123818

# Row 10
$ws.Range("I10").Value = 0.1369374790620155
$ws.Range("J10").Value = 0.1369374790620154
$ws.Range("M10").Value = 1.162287333333333
$ws.Range("N10").Value = 3.486862
$ws.Range("O10").Value = 0.2284193176890152
$ws.Range("P10").Value = 0.2284193176890152
$ws.Range("Q10").Value = 92.78220462926623
$ws.Range("R10").Value = 835.039841663396
$ws.Range("S10").Value = 0.03127916553339938
$ws.Range("T10").Value = 0.03127916553339937

# Row 11
$ws.Range("I11").Value = 0.1369374790620155
$ws.Range("J11").Value = 0.1369374790620154
$ws.Range("M11").Value = 0.232356
$ws.Range("N11").Value = 0.697068
$ws.Range("O11").Value = 0.04566392273134021
$ws.Range("P11").Value = 0.04566392273134022
$ws.Range("Q11").Value = 18.548341120616
$ws.Range("R11").Value = 166.935070085544
$ws.Range("S11").Value = 0.006253102462912393
$ws.Range("T11").Value = 0.006253102462912392

# Row 12
$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("M12").Value = 1.550905333333333
$ws.Range("N12").Value = 4.652716
$ws.Range("O12").Value = 0.3047927374587132
$ws.Range("P12").Value = 0.3047927374587132
$ws.Range("Q12").Value = 230.812582055524
$ws.Range("R12").Value = 2077.313238499716
$ws.Range("S12").Value = 0.07781260415349926
$ws.Range("T12").Value = 0.07781260415349926

# Row 13
$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("O13").Value = 0.3698654061072322
$ws.Range("P13").Value = 0.3698654061072323
$ws.Range("Q13").Value = 280.090628498618
$ws.Range("R13").Value = 2520.815656487562
$ws.Range("S13").Value = 0.09442544686417877
$ws.Range("T13").Value = 0.09442544686417878

# Row 14
$ws.Range("G14").Value = 148.824417
$ws.Range("H14").Value = 446.473251
$ws.Range("I14").Value = 0.2552967790580629
$ws.Range("J14").Value = 0.2552967790580629
$ws.Range("M14").Value = 0.260824
$ws.Range("N14").Value = 0.7824719999999999
$ws.Range("O14").Value = 0.05125861601369915
$ws.Range("P14").Value = 0.05125861601369915
$ws.Range("Q14").Value = 38.816979739608
$ws.Range("R14").Value = 349.352817656472
$ws.Range("S14").Value = 0.01308615956727143
$ws.Range("T14").Value = 0.01308615956727143

# Row 15
$ws.Range("G15").Value = 148.824417
$ws.Range("H15").Value = 446.473251
$ws.Range("I15").Value = 0.2552967790580629
$ws.Range("J15").Value = 0.2552967790580629
$ws.Range("M15").Value = 1.162287333333333
$ws.Range("N15").Value = 3.486862
$ws.Range("O15").Value = 0.2284193176890152
$ws.Range("P15").Value = 0.2284193176890152
$ws.Range("Q15").Value = 172.976734769818
$ws.Range("R15").Value = 1556.790612928362
$ws.Range("S15").Value = 0.05831471608064597
$ws.Range("T15").Value = 0.05831471608064597

# Row 16
$ws.Range("G16").Value = 148.824417
$ws.Range("H16").Value = 446.473251
$ws.Range("I16").Value = 0.2552967790580629
$ws.Range("J16").Value = 0.2552967790580629
$ws.Range("M16").Value = 0.232356
$ws.Range("N16").Value = 0.697068
$ws.Range("O16").Value = 0.04566392273134021
$ws.Range("P16").Value = 0.04566392273134022
$ws.Range("Q16").Value = 34.58024623645201
$ws.Range("R16").Value = 311.222216128068
$ws.Range("S16").Value = 0.01165785239246742
$ws.Range("T16").Value = 0.01165785239246742

# Row 17
$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("M17").Value = 1.550905333333333
$ws.Range("N17").Value = 4.652716
$ws.Range("O17").Value = 0.3047927374587132
$ws.Range("P17").Value = 0.3047927374587132
$ws.Range("Q17").Value = 54.943695260916
$ws.Range("R17").Value = 494.493257348244
$ws.Range("S17").Value = 0.01852287241879945
$ws.Range("T17").Value = 0.01852287241879945

# Row 18
$ws.Range("G18").Value = 35.426853
$ws.Range("H18").Value = 106.280559
$ws.Range("I18").Value = 0.06077202683121193
$ws.Range("J18").Value = 0.06077202683121192
$ws.Range("O18").Value = 0.3698654061072322
$ws.Range("P18").Value = 0.3698654061072323
$ws.Range("Q18").Value = 66.67406950096201
$ws.Range("R18").Value = 600.0666255086579
$ws.Range("S18").Value = 0.02247747038388581
$ws.Range("T18").Value = 0.02247747038388581

# Row 19
$ws.Range("G19").Value = 35.426853
$ws.Range("H19").Value = 106.280559
$ws.Range("I19").Value = 0.06077202683121193
$ws.Range("J19").Value = 0.06077202683121192
$ws.Range("M19").Value = 0.260824
$ws.Range("N19").Value = 0.7824719999999999
$ws.Range("O19").Value = 0.05125861601369915
$ws.Range("P19").Value = 0.05125861601369915
$ws.Range("Q19").Value = 9.240173506872001
$ws.Range("R19").Value = 83.16156156184799
$ws.Range("S19").Value = 0.003115089987715314
$ws.Range("T19").Value = 0.003115089987715314

# Row 20
$ws.Range("G20").Value = 35.426853
$ws.Range("H20").Value = 106.280559
$ws.Range("I20").Value = 0.06077202683121193
$ws.Range("J20").Value = 0.06077202683121192
$ws.Range("M20").Value = 1.162287333333333
$ws.Range("N20").Value = 3.486862
$ws.Range("O20").Value = 0.2284193176890152
$ws.Range("P20").Value = 0.2284193176890152
$ws.Range("Q20").Value = 41.17618250176201
$ws.Range("R20").Value = 370.585642515858
$ws.Range("S20").Value = 0.01388150490336395
$ws.Range("T20").Value = 0.01388150490336395

# Row 21
$ws.Range("G21").Value = 35.426853
$ws.Range("H21").Value = 106.280559
$ws.Range("I21").Value = 0.06077202683121193
$ws.Range("J21").Value = 0.06077202683121192
$ws.Range("M21").Value = 0.232356
$ws.Range("N21").Value = 0.697068
$ws.Range("O21").Value = 0.04566392273134021
$ws.Range("P21").Value = 0.04566392273134022
$ws.Range("Q21").Value = 8.231641855668
$ws.Range("R21").Value = 74.084776701012
$ws.Range("S21").Value = 0.002775089137447396
$ws.Range("T21").Value = 0.002775089137447396

# Row 22
$ws.Range("G22").Value = 121.3248153333333
$ws.Range("H22").Value = 363.974446
$ws.Range("I22").Value = 0.2081233388901116
$ws.Range("J22").Value = 0.2081233388901115
$ws.Range("M22").Value = 1.550905333333333
$ws.Range("N22").Value = 4.652716
$ws.Range("O22").Value = 0.3047927374587132
$ws.Range("P22").Value = 0.3047927374587132
$ws.Range("Q22").Value = 188.1633031661484
$ws.Range("R22").Value = 1693.469728495336
$ws.Range("S22").Value = 0.06343448218936457
$ws.Range("T22").Value = 0.06343448218936457

# Row 23
$ws.Range("G23").Value = 121.3248153333333
$ws.Range("H23").Value = 363.974446
$ws.Range("I23").Value = 0.2081233388901116
$ws.Range("J23").Value = 0.2081233388901115
$ws.Range("O23").Value = 0.3698654061072322
$ws.Range("P23").Value = 0.3698654061072323
$ws.Range("Q23").Value = 228.3358098368502
$ws.Range("R23").Value = 2055.022288531652
$ws.Range("S23").Value = 0.07697762325898423
$ws.Range("T23").Value = 0.07697762325898423

# Row 24
$ws.Range("G24").Value = 121.3248153333333
$ws.Range("H24").Value = 363.974446
$ws.Range("I24").Value = 0.2081233388901116
$ws.Range("J24").Value = 0.2081233388901115
$ws.Range("M24").Value = 0.260824
$ws.Range("N24").Value = 0.7824719999999999
$ws.Range("O24").Value = 0.05125861601369915
$ws.Range("P24").Value = 0.05125861601369915
$ws.Range("Q24").Value = 31.64442363450133
$ws.Range("R24").Value = 284.799812710512
$ws.Range("S24").Value = 0.01066811431165721
$ws.Range("T24").Value = 0.01066811431165721

# Row 25
$ws.Range("G25").Value = 121.3248153333333
$ws.Range("H25").Value = 363.974446
$ws.Range("I25").Value = 0.2081233388901116
$ws.Range("J25").Value = 0.2081233388901115
$ws.Range("M25").Value = 1.162287333333333
$ws.Range("N25").Value = 3.486862
$ws.Range("O25").Value = 0.2284193176890152
$ws.Range("P25").Value = 0.2284193176890152
$ws.Range("Q25").Value = 141.0142960809391
$ws.Range("R25").Value = 1269.128664728452
$ws.Range("S25").Value = 0.04753939106443895
$ws.Range("T25").Value = 0.04753939106443895

# Row 26
$ws.Range("G26").Value = 121.3248153333333
$ws.Range("H26").Value = 363.974446
$ws.Range("I26").Value = 0.2081233388901116
$ws.Range("J26").Value = 0.2081233388901115
$ws.Range("M26").Value = 0.232356
$ws.Range("N26").Value = 0.697068
$ws.Range("O26").Value = 0.04566392273134021
$ws.Range("P26").Value = 0.04566392273134022
$ws.Range("Q26").Value = 28.190548791592
$ws.Range("R26").Value = 253.714939124328
$ws.Range("S26").Value = 0.009503728065666588
$ws.Range("T26").Value = 0.009503728065666588
